$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update score values
$ws.Range("D2").Value = 50
$ws.Range("F2").Value = 258

$ws.Range("C3").Value = 96
$ws.Range("D3").Value = 25
$ws.Range("F3").Value = 229

$ws.Range("C4").Value = 80
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = 197

$ws.Range("C5").Value = 56
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 185

$ws.Range("B6").Value = 65
$ws.Range("D6").Value = 30
$ws.Range("F6").Value = 183

# Update selected cell
$ws.Range("F7").Select()
